$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.682.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.556.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("E6").Value = "  -2.06%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.94"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.778.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.563.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("E14").Value = "  -2.33%  "
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.734.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "213.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0674"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("E22").Value = "  -0.85%  "
$ws.Range("E23").Value = "  -1.73%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("E26").Value = "  +0.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("E29").Value = "  -0.99%  "
$ws.Range("E30").Value = "  -1.62%  "
$ws.Range("E31").Value = "  -3.89%  "
$ws.Range("E32").Value = "  -1.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.384.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.49%  "
$ws.Range("E34").Value = "  -1.43%  "
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("E37").Value = "  -4.21%  "
$ws.Range("E38").Value = "  -2.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.519"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.810"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("E42").Value = "  +2.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.87%  "
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("E45").Value = "  -1.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.691.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("E49").Value = "  -2.11%  "
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0948"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.87%  "
